$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Vi): shift existing values two columns right (from E2 onward) and
# insert two new data points (0.6 and 0.82), extending the series out to X2.
$ws.Range("B2").Value = 0
$ws.Range("C2").Value = 0.25
$ws.Range("D2").Value = 0.5
$ws.Range("E2").Value = 0.6
$ws.Range("F2").Value = 0.75
$ws.Range("G2").Value = 0.82
$ws.Range("H2").Value = 1
$ws.Range("I2").Value = 1.1
$ws.Range("J2").Value = 1.2
$ws.Range("K2").Value = 1.3
$ws.Range("L2").Value = 1.4
$ws.Range("M2").Value = 1.5
$ws.Range("N2").Value = 1.55
$ws.Range("O2").Value = 1.6
$ws.Range("P2").Value = 1.65
$ws.Range("Q2").Value = 1.7
$ws.Range("R2").Value = 1.75
$ws.Range("S2").Value = 1.8
$ws.Range("T2").Value = 1.9
$ws.Range("U2").Value = 2
$ws.Range("V2").Value = 3
$ws.Range("W2").Value = 4
$ws.Range("X2").Value = 5

# Row 3 (Vo): new measured values, including filling in B3:E3 which were
# previously empty, and shifting the rest two columns right.
$ws.Range("B3").Value = 4.91
$ws.Range("C3").Value = 4.91
$ws.Range("D3").Value = 4.91
$ws.Range("E3").Value = 4.78
$ws.Range("F3").Value = 4.22
$ws.Range("G3").Value = 3.85
$ws.Range("H3").Value = 2.25
$ws.Range("I3").Value = 2.03
$ws.Range("J3").Value = 1.47
$ws.Range("K3").Value = 0.95
$ws.Range("L3").Value = 0.54
$ws.Range("M3").Value = 0.24
$ws.Range("N3").Value = 0.32
$ws.Range("O3").Value = 0.24
$ws.Range("P3").Value = 0.21
$ws.Range("Q3").Value = 0.2
$ws.Range("R3").Value = 0.191
$ws.Range("S3").Value = 0.181
$ws.Range("T3").Value = 1.171
$ws.Range("U3").Value = 0.163
$ws.Range("V3").Value = 0.123
$ws.Range("W3").Value = 0.106
$ws.Range("X3").Value = 0.095

# New section header for a second (PNP) data table below, on row 5.
$ws.Range("A5").Value = "PNP"

# Update the view: drop the frozen/scrolled top-left cell and move the
# active selection.
$ws.Range("G8").Select()
